# Restored from revision #195ae0bf3a96a4be83d88e277669f05aa2008e54.TEST
# Author: admin. Type: SAVE.
#
# The "Rules" sheet drives a lookup table (Rule / From / To / Greeting).
# Row 10 is rule "R30": its "From" threshold (cell C10) is corrected
# from 18 back to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
